# Adding new progress as of date 04 Nov 2025
# For each data row (3-29) on the Training Dashboard sheet:
#   - decrement the "PERIOD TO EXPIRE" (column H) value by 1
#   - bump the "LAST UPDATE" (column I) date text from 03-Nov-2025 to 04-Nov-2025
#     (kept as literal text, matching how the source data is stored)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 29; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $current = $hCell.Value2()
    $hCell.Value = $current - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}
